# Updates cryptos list values (price / 1h volume change) per the
# "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.005.77'
$ws.Range("E2").Value = '  -1.86%  '
# Row 3
$ws.Range("D3").Value = '1.977.10'
$ws.Range("E3").Value = '  -3.65%  '
# Row 4
$ws.Range("E4").Value = '  +0.18%  '
# Row 5
$ws.Range("D5").Value = '''242.58'
$ws.Range("E5").Value = '  -8.54%  '
# Row 6
$ws.Range("D6").Value = '''0.599'
$ws.Range("E6").Value = '  -4.44%  '
# Row 7
$ws.Range("E7").Value = '  +0.17%  '
# Row 8
$ws.Range("D8").Value = '''54.38'
$ws.Range("E8").Value = '  -7.46%  '
# Row 9
$ws.Range("D9").Value = '''0.370'
$ws.Range("E9").Value = '  -5.55%  '
# Row 10
$ws.Range("D10").Value = '''0.0748'
$ws.Range("E10").Value = '  -8.08%  '
# Row 11
$ws.Range("D11").Value = '''0.0992'
$ws.Range("E11").Value = '  -4.40%  '
# Row 12
$ws.Range("D12").Value = '2.267.68'
$ws.Range("E12").Value = '  -3.00%  '
# Row 13
$ws.Range("E13").Value = '  -8.50%  '
# Row 14
$ws.Range("D14").Value = '''20.89'
$ws.Range("E14").Value = '  -4.00%  '
# Row 15
$ws.Range("D15").Value = '''0.756'
$ws.Range("E15").Value = '  -9.99%  '
# Row 16
$ws.Range("D16").Value = '''5.06'
$ws.Range("E16").Value = '  -7.25%  '
# Row 17
$ws.Range("D17").Value = '1.966.54'
$ws.Range("E17").Value = '  -4.15%  '
# Row 18
$ws.Range("D18").Value = '36.836.31'
$ws.Range("E18").Value = '  -2.05%  '
# Row 19
$ws.Range("D19").Value = '''68.44'
$ws.Range("E19").Value = '  -3.09%  '
# Row 20
$ws.Range("E20").Value = '  -6.36%  '
# Row 21
$ws.Range("D21").Value = '''227.87'
$ws.Range("E21").Value = '  -1.06%  '
# Row 22
$ws.Range("D22").Value = '''4.96'
$ws.Range("E22").Value = '  -6.80%  '
# Row 23
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  -0.11%  '
# Row 24
$ws.Range("D24").Value = '''2.42'
$ws.Range("E24").Value = '  -9.36%  '
# Row 25
$ws.Range("D25").Value = '''2.34'
$ws.Range("E25").Value = '  -0.87%  '
# Row 26
$ws.Range("D26").Value = '''162.49'
$ws.Range("E26").Value = '  -1.17%  '
# Row 27
$ws.Range("D27").Value = '''8.63'
$ws.Range("E27").Value = '  -7.61%  '
# Row 28
$ws.Range("E28").Value = '  -8.80%  '
# Row 29
$ws.Range("D29").Value = '''19.09'
$ws.Range("E29").Value = '  -4.91%  '
# Row 30
$ws.Range("E30").Value = '  -7.41%  '
# Row 31
$ws.Range("E31").Value = '  -3.98%  '
# Row 32
$ws.Range("E32").Value = '  -8.27%  '
# Row 33
$ws.Range("D33").Value = '''0.0611'
$ws.Range("E33").Value = '  -9.58%  '
# Row 34
$ws.Range("D34").Value = '''4.28'
$ws.Range("E34").Value = '  -6.58%  '
# Row 35
$ws.Range("D35").Value = '''2.35'
$ws.Range("E35").Value = '  -6.86%  '
# Row 36
$ws.Range("E36").Value = '  -0.12%  '
# Row 37
$ws.Range("E37").Value = '  +0.30%  '
# Row 38
$ws.Range("D38").Value = '''3.33'
$ws.Range("E38").Value = '  -6.82%  '
# Row 39
$ws.Range("D39").Value = '''5.16'
$ws.Range("E39").Value = '  -5.45%  '
# Row 40
$ws.Range("E40").Value = '  +0.19%  '
# Row 41
$ws.Range("D41").Value = '1.420.65'
$ws.Range("E41").Value = '  +0.87%  '
# Row 42
$ws.Range("E42").Value = '  -4.91%  '
# Row 43
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").Value = '''0.0898'
$ws.Range("E43").Value = '  -8.65%  '
# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '''0.0204'
$ws.Range("E44").Value = '  -7.24%  '
# Row 45
$ws.Range("D45").Value = '''15.50'
$ws.Range("E45").Value = '  -7.84%  '
# Row 46
$ws.Range("D46").Value = '''87.35'
$ws.Range("E46").Value = '  -5.60%  '
# Row 47
$ws.Range("D47").Value = '''0.998'
$ws.Range("E47").Value = '  -6.56%  '
# Row 48
$ws.Range("D48").Value = '''2.87'
$ws.Range("E48").Value = '  -1.18%  '
# Row 49
$ws.Range("D49").Value = '''6.63'
$ws.Range("E49").Value = '  -12.03%  '
# Row 50
$ws.Range("D50").Value = '2.162.46'
$ws.Range("E50").Value = '  -2.97%  '
# Row 51
$ws.Range("D51").Value = '''1.85'
$ws.Range("E51").Value = '  -12.05%  '
